# [PHOENIX-6108] Completed the Bank to Bank payment transfer scenario
#
# legalCaseTestData.xlsx: the petitionType test value on row 2 changes
# from "EXECUTION PETITION" to "CIVIL MISCELLANEOUS PETITION", the sheet's
# active selection moves to D2, and the column widths are widened to
# accommodate the longer text (most visibly column C / petitionType).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the petitionType value for testData1 (cell C2).
$ws.Range("C2").Value = "CIVIL MISCELLANEOUS PETITION"

# New active cell / selection.
$ws.Range("D2").Select()

# Resulting (auto-adjusted) column widths, column C widening the most
# since it now holds the much longer "CIVIL MISCELLANEOUS PETITION" text.
$ws.Columns.Item(1).ColumnWidth = 11.414814814814767
$ws.Columns.Item(2).ColumnWidth = 12.985185185185166
$ws.Columns.Item(3).ColumnWidth = 29.02222222222227
$ws.Columns.Item(4).ColumnWidth = 20.529629629629667
$ws.Columns.Item(5).ColumnWidth = 18.57037037037037
$ws.Columns.Item(6).ColumnWidth = 15.629629629629667
$ws.Columns.Item(7).ColumnWidth = 14.844444444444466
$ws.Columns.Item(8).ColumnWidth = 16.41481481481477
$ws.Columns.Item(9).ColumnWidth = 17.981481481481467
$ws.Columns.Item(10).ColumnWidth = 21.80370370370367
$ws.Columns.Item(11).ColumnWidth = 17.49259259259257
$ws.Columns.Item(12).ColumnWidth = 17.588888888888867

# Remaining columns' default width also grows slightly.
$ws.StandardWidth = 8.377777777777776
